$d = $word.ActiveDocument

$old = "Allow add XY for lines and polygons as well"
$new = "Allow add XY for lines and polygons as well.^pIf there are a lot of points, provide a progress bar^pIf data is in a compressed FGDB, it is read-only, therefore, it will not show up in the list of available features (this can be confusing)."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
